$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column E
$ws.Range("E1").Value = "aggregate_id"

# Fill column E (rows 2-15) with the same year values as column A
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 1).Value2
}

# Update the active selection to match the target state
$ws.Range("G12").Select()
